$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.283.92'
$ws.Range('E2').Value = '  +5.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.801.23'
$ws.Range('E3').Value = '  +3.49%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.48'
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5512'
$ws.Range('E7').Value = '  +9.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3867'
$ws.Range('E8').Value = '  +8.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07564'
$ws.Range('E9').Value = '  +4.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.73'
$ws.Range('E10').Value = '  +1.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.119'
$ws.Range('E11').Value = '  +5.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.0000'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  +5.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.199'
$ws.Range('E14').Value = '  +4.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.349'
$ws.Range('E15').Value = '  +8.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.799.03'
$ws.Range('E16').Value = '  +3.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.89'
$ws.Range('E17').Value = '  +6.21%  '
$ws.Range('E18').Value = '  +3.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06449'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.27'
$ws.Range('E21').Value = '  +4.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.973'
$ws.Range('E22').Value = '  +4.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.307.33'
$ws.Range('E23').Value = '  +5.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.44'
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.146'
$ws.Range('E25').Value = '  +5.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.86'
$ws.Range('E26').Value = '  +3.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.65'
$ws.Range('E27').Value = '  +4.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.399'
$ws.Range('E28').Value = '  +8.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.007.10'
$ws.Range('E29').Value = '  +3.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.50'
$ws.Range('E30').Value = '  +3.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.124'
$ws.Range('E31').Value = '  +8.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1021'
$ws.Range('E32').Value = '  +6.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.722'
$ws.Range('E33').Value = '  +6.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.672'
$ws.Range('E34').Value = '  +2.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.2373'
$ws.Range('E35').Value = '  +18.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06329'
$ws.Range('E36').Value = '  +7.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.885'
$ws.Range('E37').Value = '  +17.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02317'
$ws.Range('E38').Value = '  +6.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.61'
$ws.Range('E39').Value = '  +5.86%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.031'
$ws.Range('E40').Value = '  +5.87%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6383'
$ws.Range('E41').Value = '  +5.82%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.156'
$ws.Range('E43').Value = '  +4.56%  '
$ws.Range('E44').Value = '  -2.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.49'
$ws.Range('E45').Value = '  +5.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5965'
$ws.Range('E46').Value = '  +5.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.683'
$ws.Range('E47').Value = '  +2.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.74'
$ws.Range('E48').Value = '  +3.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.976'
$ws.Range('E49').Value = '  +7.33%  '
$ws.Range('E50').Value = '  +4.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06894'
$ws.Range('E51').Value = '  +3.61%  '
